$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.498.36'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.356.24'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '256.42'
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '664.75'
$ws.Range("E6").Value = '  +6.41%  '
$ws.Range("E7").Value = '  +6.20%  '
$ws.Range("E8").Value = '  +20.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.09'
$ws.Range("E9").Value = '  +21.50%  '
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.353.02'
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("E12").Value = '  +8.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.35'
$ws.Range("E13").Value = '  +12.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000271'
$ws.Range("E14").Value = '  +8.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '97.901.51'
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.979.46'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.66'
$ws.Range("E17").Value = '  +2.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.89'
$ws.Range("E18").Value = '  +28.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.342.38'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.94'
$ws.Range("E20").Value = '  +11.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '529.38'
$ws.Range("E21").Value = '  +7.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.54'
$ws.Range("E22").Value = '  -0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.49'
$ws.Range("E23").Value = '  +11.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000215'
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.442'
$ws.Range("E25").Value = '  +47.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '102.83'
$ws.Range("E26").Value = '  +11.96%  '
$ws.Range("E27").Value = '  +10.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.57'
$ws.Range("E28").Value = '  +5.19%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.538.35'
$ws.Range("E29").Value = '  +0.58%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.150'
$ws.Range("E30").Value = '  +9.36%  '
$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.21'
$ws.Range("E32").Value = '  +15.00%  '
$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.190'
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '29.49'
$ws.Range("E35").Value = '  +4.43%  '
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.544'
$ws.Range("E36").Value = '  +17.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.12'
$ws.Range("E37").Value = '  +8.45%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.74'
$ws.Range("E38").Value = '  +6.43%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.157'
$ws.Range("E39").Value = '  +5.16%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '521.81'
$ws.Range("E40").Value = '  +3.78%  '
$ws.Range("B41").Value = 'MantraDAO'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.93'
$ws.Range("E41").Value = '  +7.35%  '
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.34'
$ws.Range("E42").Value = '  +5.93%  '
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '24.72'
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0436'
$ws.Range("E44").Value = '  +33.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.828'
$ws.Range("E45").Value = '  +6.34%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.39'
$ws.Range("E46").Value = '  +1.70%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.18'
$ws.Range("E48").Value = '  +11.00%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.07'
$ws.Range("E49").Value = '  +6.20%  '
$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.55'
$ws.Range("E50").Value = '  +13.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.88'
$ws.Range("E51").Value = '  +14.86%  '
